# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 23:07"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4302867
$ws.Range("C4").Value = 54540
$ws.Range("D4").Value = 2052197
$ws.Range("E4").Value = 2101436
$ws.Range("G4").Value = 744
$ws.Range("H4").Value = 149234

# Row 8 - Sudafrica
$ws.Range("B8").Value = 434200
$ws.Range("C8").Value = 12204
$ws.Range("D8").Value = 263054
$ws.Range("E8").Value = 164491
$ws.Range("G8").Value = 312
$ws.Range("H8").Value = 6655

# Row 21 - Alemania
$ws.Range("B21").Value = 206332
$ws.Range("C21").Value = 372
$ws.Range("E21").Value = 6730

# Row 28 - Egipto
$ws.Range("B28").Value = 91583
$ws.Range("C28").Value = 511
$ws.Range("D28").Value = 32903
$ws.Range("E28").Value = 54122
$ws.Range("G28").Value = 40
$ws.Range("H28").Value = 4558

# Row 41 - Israel
$ws.Range("B41").Value = 60678
$ws.Range("C41").Value = 1203
$ws.Range("D41").Value = 26917
$ws.Range("E41").Value = 33304
$ws.Range("G41").Value = 9
$ws.Range("H41").Value = 457

# Row 51 - Barein
$ws.Range("B51").Value = 38747
$ws.Range("C51").Value = 289
$ws.Range("D51").Value = 35205
$ws.Range("E51").Value = 3405

# Row 80 - Estado de Palestina
$ws.Range("E80").Value = 6949
$ws.Range("G80").Value = 5
$ws.Range("H80").Value = 75

# Row 134 - Yemen
$ws.Range("D134").Value = 780
$ws.Range("E134").Value = 420
$ws.Range("G134").Value = 5
$ws.Range("H134").Value = 474

$wb.Save()
